$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 3's date/link text (E3, F3) ---
$ws.Cells.Item(3, 5).Value = "18/25/2022 00:25:09"
$ws.Cells.Item(3, 6).Value = "https://www.plus2net.com"

# --- Add new row 8 ---
$ws.Cells.Item(8, 1).Value = 8

# B8 must stay a text value ("123345") rather than be auto-coerced to a number
$bCell = $ws.Cells.Item(8, 2)
$bCell.NumberFormat = "@"
$bCell.Value = "123345"
$bCell.Style = "Normal"

$ws.Cells.Item(8, 3).Value = "<p>1235</p>"
$ws.Cells.Item(8, 4).Value = "Ban Đào Tạo"
$ws.Cells.Item(8, 5).Value = "17/51/2022 22:51:28"
$ws.Cells.Item(8, 6).Value = "https://www.plus2net.com"

# --- Update the active selection to match the new edit location ---
$ws.Range("E9").Select() | Out-Null
